# "Gebruik DPIA in plaats van PIA." (#1024)
# Replace the "PIA" abbreviation with "DPIA" everywhere it appears on the
# single slide of this deck, and refresh the cached datetimeFigureOut
# fields on the Handout Master / Notes Master to the new save date.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1) "Ezelsoor 4" folded-corner shape: PIA -> DPIA -----------------
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Ezelsoor 4") {
        $sh.TextFrame.TextRange.Text = "DPIA"
    }
}

# --- 2) "Tekstvak 109" legend textbox: expand the PIA explanation -----
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Tekstvak 109") {
        $tr = $sh.TextFrame.TextRange

        # Locate the paragraph that still reads "PIA<TAB>privacy impact analyse"
        for ($pi = 1; $pi -le $tr.Paragraphs().Count; $pi++) {
            $para = $tr.Paragraphs($pi, 1)
            if ($para.Text -like "PIA`t*") {
                $paraStart = $para.Start

                # Rewrite the whole paragraph text first (this keeps the
                # common " impact analyse" suffix in its own run and puts
                # everything else into a single new run).
                $para.Text = "DPIA`tdata protection impact analyse"

                # Re-split "DPIA`tdata protection" so "protection" becomes
                # its own run (mirrors the err="1" spellcheck run in the
                # authored deck), by re-assigning just that sub-range.
                $word = $tr.Characters($paraStart + 10, 10)
                $word.Text = "protection"
            }
        }
    }
}

# --- 3) Handout Master / Notes Master date placeholders ---------------
# Re-display the "datetimeFigureOut" field with the current save date.
$newDate = "24-03-2025"

$handoutMaster = $p.HandoutMaster
$handoutMaster.HeadersFooters.DateAndTime.Text = $newDate

$notesMaster = $p.NotesMaster
$notesMaster.HeadersFooters.DateAndTime.Text = $newDate
